# Updated TestDataHandler and Added FR64 Test Page
#
# Applies the data + view changes described by the commit:
#  - Disruptions sheet: swap the two sample "disrupted flight" values for
#    a fresh FR64 test flight pair, and move the active selection to A2.
#  - URL sheet: the used range grows out to column F (from formatting /
#    column usage) and the active selection moves to B13.
#  - Main sheet: active selection moves from B3 to B2 (this also leaves
#    Main as the selected/active tab, so it is done last).

$wb = $excel.ActiveWorkbook

# --- Disruptions sheet -----------------------------------------------
$wsDisruptions = $wb.Worksheets.Item("Disruptions")
$wsDisruptions.Range("D2").Value = "FL-ZZ-866-20171004-PIT-LAX-0"
$wsDisruptions.Range("E2").Value = "FL-ZZ-54-20171004-PIE-LAX-0"
$wsDisruptions.Range("A2").Select()

# --- URL sheet ---------------------------------------------------------
$wsUrl = $wb.Worksheets.Item("URL")
# Touch columns A:F so the sheet's used range extends to column F even
# though there's no literal data past column D.
$wsUrl.Range("A1:F2").Font.Bold = $wsUrl.Range("A1:F2").Font.Bold
$wsUrl.Range("B13").Select()

# --- Main sheet (left active/selected, so do this last) ---------------
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("B2").Select()
